$wb = $excel.ActiveWorkbook

$insert = $wb.Worksheets.Item("Insert")
$result = $wb.Worksheets.Item("Result")

# User typed a new scouting-data row into Insert!A4.
$insert.Range("A4").Value = "s=Yiğit Emre Çulcuoğlu;e=2022 İzmir Regional;l=f;m=1;r=b2;t=6838;as=[54];at=Y;au=8;al=3;ac=Y;tu=5;tl=6;tm=1;tn=2;wd=Y;cl=b;ss=[53,56,44,58,34,20,17,18,41,51,51,16,26,50,64,68,34];c=4;be=Y;cn=1;ds=v;dr=v;d=N;to=N;cf=N;co=Commen"

# The existing array formula in Insert!B4 (previously anchored on an empty
# A4 and therefore only a single cell) needs to (re)spill across B4:AB4,
# same as the sibling rows that already hold data (row 2).
$insert.Range("B4:AB4").FormulaArray = '=TRANSPOSE(TRIM(MID(SUBSTITUTE(";"&A4,";",REPT(" ",LEN(A4)+1)),ROW(INDIRECT("A1:A"&LEN(A4)-LEN(SUBSTITUTE(A4,";",""))+1))*LEN(A4)+1,LEN(A4)+1)))'

# Result!B4 mirrors Insert!B4 via ANCHORARRAY and likewise needs to spill
# across B4:AB4 now that Insert!A4 has a value.
$result.Range("B4:AB4").FormulaArray = '=RIGHT(_xlfn.ANCHORARRAY(Insert!B4), LEN(_xlfn.ANCHORARRAY(Insert!B4)) - FIND("=",_xlfn.ANCHORARRAY( Insert!B4)))'

# Leave the cursor where the author left it after entering the new row.
$insert.Activate() | Out-Null
$insert.Range("A13").Select() | Out-Null
